# Adds the "Initial" column to the Node_start sheet (the new decision-variable
# column the commit message refers to: "añadí la columna de inicial y agregué
# las variables de decisión") and keeps the sheet's filter / defined-name
# metadata in sync with the now-wider data range.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Node_start: add column E = "Initial" with its single data value 1000
# ---------------------------------------------------------------------
$wsStart = $wb.Worksheets.Item("Node_start")

# Match the look of the existing table: header cells (row 1) share the bold
# white-on-blue style, body cells (row 2) share the light banded-fill style.
$wsStart.Range("C1").Copy()
$wsStart.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$wsStart.Range("A2").Copy()
$wsStart.Range("E2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$wsStart.Range("E1").Value = "Initial"
$wsStart.Range("E2").Value = 1000

# The sheet's AutoFilter + hidden _FilterDatabase name must grow from
# A1:D2 to A1:E2 to cover the new column.
$wsStart.AutoFilterMode = $false
$wsStart.Range("A1:E2").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Node_start!_FilterDatabase") {
        $n.RefersTo = "=Node_start!`$A`$1:`$E`$2"
    }
}

# ---------------------------------------------------------------------
# 2. Leave the workbook the way the author left it: on the Arcs sheet.
# ---------------------------------------------------------------------
$wsStart.Range("F10").Select()

$wsNodeEnd = $wb.Worksheets.Item("Node_end")
$wsNodeEnd.Activate()
$wsNodeEnd.Range("G10").Select()

$wsNodes = $wb.Worksheets.Item("Nodes")
$wsNodes.Activate()
$wsNodes.Range("H8").Select()

$wsArcs = $wb.Worksheets.Item("Arcs")
$wsArcs.Activate()
$wsArcs.Range("M8").Select()
